$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest three years (2007-2009), shifting all rows up by 3.
$ws.Rows("2:4").Delete()

# Append the newly reported year (2021) as the new last row (row 13).
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 48
$ws.Range("C13").Value = 574
$ws.Range("D13").Value = 4006
$ws.Range("F13").Value = 5319
$ws.Range("H13").Value = 2703
$ws.Range("K13").Value = 428
$ws.Range("L13").Value = 179

# Match the formatting (bold, bordered, centered) used by the other year labels.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

Write-Output "A13=$($ws.Range('A13').Value())"
Write-Output "dim=$($ws.UsedRange.Address())"
